# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Updates the "Periodo Mora" (E), "Valor Mora" (F) and "Salario Basico" (G)
# columns for the worker's debt-period rows (16-22) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Periodo Mora" labels, now listed in ascending chronological order.
$periodos = @("2008", "2010", "2101", "2102", "2106", "2107", "2108")

# New "Valor Mora" amounts lined up with the re-ordered periods above.
$valores = @(35112, 35112, 35112, 35112, 36341, 36341, 33942)

# New uniform "Salario Basico" amount for every period row.
$salario = 877803

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
    $ws.Cells.Item($row, 7).Value = $salario
}
